$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sun Dec 17 18_34_08 2023"
$ws.Range("B2").Value = "loc"
$ws.Range("C2").Value = 180

$ws.Range("A3").Value = "Sun Dec 17 19_23_37 2023"
$ws.Range("B3").Value = "loc2"
$ws.Range("C3").Value = -30
